# Update "想去人数" (want-to-go count) values in column F across the
# workbook's sheets, matching the regenerated data pulled at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# sheet name -> list of (row, newValue) updates for column F
$updates = @{
    "展览"     = @(
        @{ Row = 11; Value = 663 },
        @{ Row = 13; Value = 5863 },
        @{ Row = 17; Value = 5564 },
        @{ Row = 25; Value = 96 },
        @{ Row = 26; Value = 1166 },
        @{ Row = 31; Value = 3833 }
    )
    "演出"     = @(
        @{ Row = 5; Value = 201 }
    )
    "本地生活" = @(
        @{ Row = 2; Value = 9455 },
        @{ Row = 4; Value = 2182 }
    )
    "全部类型" = @(
        @{ Row = 2; Value = 9455 },
        @{ Row = 4; Value = 2182 },
        @{ Row = 14; Value = 663 },
        @{ Row = 16; Value = 5863 },
        @{ Row = 25; Value = 5564 },
        @{ Row = 33; Value = 96 },
        @{ Row = 34; Value = 1166 },
        @{ Row = 45; Value = 3833 }
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($entry in $updates[$sheetName]) {
        $ws.Range("F" + $entry.Row).Value = $entry.Value
    }
}

$wb.Save()
